$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("marker_table")

# Collapse homozygous genotype strings (e.g. "GG") down to a single
# allele letter (e.g. "G") in the genotype column (G) of the marker table.
$updates = @{
    "G4"  = "T"
    "G5"  = "G"
    "G6"  = "G"
    "G7"  = "G"
    "G8"  = "G"
    "G9"  = "C"
    "G10" = "G"
    "G12" = "G"
    "G13" = "A"
    "G14" = "C"
    "G15" = "A"
    "G16" = "T"
    "G17" = "G"
    "G18" = "G"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
